$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" — copy formatting (bold/border/alignment, style index)
# from the neighboring header cell G1 so it matches the existing header row style.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New data cell H2 — numeric value 1 (no special style, matches the other data cells)
$ws.Range("H2").Value = 1
